# Update "detection field test data.xlsx" per the target diff.
# - dog sheet (ActiveSheet): C8 becomes a text note, two new rows (9 & 10)
#   are appended with field-test data, column G is widened, and the
#   view/selection is reset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C8: time value replaced by a text note, keep its existing (time) style ---
$ws.Range("C8").Value = "previous morning"

# --- Row 9: new field-test entry ---
# Copy formatting from row 4 first (A=date style, C/D=time style, rest plain)
$ws.Range("A4:L4").Copy() | Out-Null
$ws.Range("A9:L9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
# C9 holds text (no number format) like C7 -- clear the copied time style on C9
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$ws.Range("A9").Value = 45796
$ws.Range("B9").Value = "PRESENCE"
$ws.Range("C9").Value = "previous evening"
$ws.Range("D9").Value = 0.45833333333333331
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 15
$ws.Range("G9").Value = "Sunny, cool"
$ws.Range("H9").Value = $false
$ws.Range("I9").Value = "16 minutes 47 seconds"
$ws.Range("J9").Value = 1007
$ws.Range("K9").Value = "NA"
$ws.Range("L9").Value = "Worked downhill. Did not get onto odour. Search got derailed by Koda picking up odour of a frisbee next to the search area. Probably going to exclude because ignoring a high value item is not part of the controlled evaluation protocol/"

# --- Row 10: new field-test entry ---
$ws.Range("A4:L4").Copy() | Out-Null
$ws.Range("A10:L10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A10").Value = 45797
$ws.Range("B10").Value = "PRESENCE"
$ws.Range("C10").Value = 0.29166666666666669
$ws.Range("D10").Value = 0.63194444444444442
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = "Sunny, mild"
$ws.Range("H10").Value = $true
$ws.Range("I10").Value = "1 minute 31 seconds"
$ws.Range("J10").Value = 91
$ws.Range("K10").Value = "Primary sweeps"
$ws.Range("L10").Value = "Worked uphill. Found on the first sweep and I saw it the moment Koda dropped into an alert."

# --- widen column G (Conditions) to fit the new content ---
# (target stored width is 13.08984375; this runtime quantises column widths
# to 1/6-character pixel steps, so 12.3 is the closest achievable setting,
# landing on a stored width of 13.1666... )
$ws.Columns.Item(7).ColumnWidth = 12.3

# --- reset the view: scroll back to the top-left and select A11 ---
$ws.Range("A11").Select() | Out-Null

Write-Host "done"
